$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card21")

$ws.Range("K7").Value = "nan"
$ws.Range("K8").Value = "✅"
